$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.5269039735099338
$ws.Range("F2").Value = 0.8436050364479788
$ws.Range("G2").Value = 0.6486624203821657

$ws.Range("E3").Value = 0.6153381642512077
$ws.Range("F3").Value = 0.6752816434724983
$ws.Range("G3").Value = 0.6439178515007898

$ws.Range("E4").Value = 0.5404178019981835
$ws.Range("F4").Value = 0.7886017229953611
$ws.Range("G4").Value = 0.6413365669630827
